$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply an AutoFilter on the data range A1:M191, filtering column F (2030, field 6)
# to show only values less than 1. This hides all rows that don't match the
# criteria (setting row hidden="1") and writes the <autoFilter> element.
$rng = $ws.Range("A1:M191")
$rng.AutoFilter(6, "<1", [Microsoft.Office.Interop.Excel.XlAutoFilterOperator]::xlAnd)

# Mark the format-conditions calculation as disabled, matching sheetPr.
$ws.EnableFormatConditionsCalculation = $false

# Register the hidden _FilterDatabase defined name used by the AutoFilter.
$ws.Names.Add("_xlnm._FilterDatabase", "=data!`$A`$1:`$M`$191", $false)
$fdName = $wb.Names.Item($wb.Names.Count)
$fdName.Visible = $false

# Update the active selection to B93.
$ws.Range("B93").Select()

# Adjust the workbook window position/size recorded in the book view.
$win = $excel.ActiveWindow
$win.Left = 240
$win.Top = 660
$win.Width = 46360
$win.Height = 25920
